$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42
$ws.Range("A42").Value = "EEEF"
$ws.Range("B42").Value = "FFDF"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "11223344"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11223344"
$ws.Range("E42").Value = "efe"
$ws.Range("F42").Value = "الماشية"
$ws.Range("G42").Value = "الأبقار"

# Row 43
$ws.Range("A43").Value = "TJT"
$ws.Range("B43").Value = "FDDF"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "11223344"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11223344"
$ws.Range("E43").Value = "DDZ"
$ws.Range("F43").Value = "الماشية"
$ws.Range("G43").Value = "الأغنام"
